# "added the end slide"
# The deck's final slide (slide 15) is a blank "end" slide with no
# shapes on it. This adds a simple autosized text box reading
# "Any Question" to that slide, matching how PowerPoint records a
# manually-inserted, auto-fit text box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)

# Target EMU box: off (3131840, 1988840) ext (1445780 x 369332).
# Shapes.AddTextbox takes Left/Top/Width/Height in points (1 pt = 12700 EMU).
$left   = 3131840 / 12700
$top    = 1988840 / 12700
$width  = 1445780 / 12700
$height = 369332 / 12700

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = "TextBox 1"

# No fill, single line (no wrap), shape grows/shrinks to fit its text.
$shp.Fill.Visible = $false
$shp.TextFrame.WordWrap = $false
$shp.TextFrame.AutoSize = 1

$shp.TextFrame.TextRange.Text = "Any Question"
$shp.TextFrame.TextRange.LanguageID = "en-IN"
